$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of H7:H80 (attendance marks for this date column).
# The daily-total formulas in column J (=SUM(E:I)) recalc automatically.
$ws.Range("H7:H80").ClearContents()
